$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (first) paragraph.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
[void]$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
[void]$metaPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Discover Book of Time, a high volatility slot game with exciting bonuses and a maximum win potential of 10,000x the bet amount. Play for free now.</w:t></w:r></w:p>')

# ---------------------------------------------------------------------------
# 2. Remove the duplicated "Play Book of Time Free: Slot Game Review"
#    paragraph that used to sit near the end of the document (the real
#    document title in paragraph 1 keeps its Heading 1 style, so skip it
#    and match the later plain-Normal duplicate instead).
# ---------------------------------------------------------------------------
$target = $null
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text.TrimEnd()
    if ($txt -eq "Play Book of Time Free: Slot Game Review") {
        $target = $para
        break
    }
}
[void]$target.Range.Delete()

# ---------------------------------------------------------------------------
# 3. Replace the text of the closing italic paragraph (the old meta
#    description) with the new image-generation prompt, keeping the
#    paragraph's italic run formatting.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
[void]$lastPara.Range.Delete()

$newLastAnchor = $d.Paragraphs.Item($d.Paragraphs.Count)
[void]$newLastAnchor.Range.InsertParagraphAfter()
$newLastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
[void]$newLastPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Prompt: Create a feature image for "Book of Time" that depicts a happy Maya warrior with glasses in a cartoon style. The image should have a colorful and vibrant background highlighting the forest theme of the game. The warrior should be holding a book and a clock symbol that represents the Scatter symbols of the game. The image should convey the retro taste and comic book style of the slot game. Use warm and inviting colors along with bold outlines to give a lively and energetic feel to the image. The image should be able to catch the player''s attention and convey the exciting and thrilling experience of playing "Book of Time".</w:t></w:r></w:p>')

Write-Host "Paragraphs: $($d.Paragraphs.Count) (expected 49)"
